# add JSON data source: insert a CustomerID column (B) into the OrderHeader
# sheet, sourced from the federated JSON customer data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Status" column (B) - this shifts
# Status -> C and Description -> D, matching the target layout.
$ws.Columns("B:B").Insert()

$ws.Range("B1").Value = "CustomerID"

$customerIds = @(
    "ANATR",
    "TORTU",
    "VAFFE",
    "MEREP",
    "PRINI",
    "HILAA",
    "LEHMS",
    "QUEEN",
    "LAMAI",
    "FRANK",
    "WARTH",
    "BONAP",
    "PERIC",
    "MEREP",
    "QUICK"
)

for ($i = 0; $i -lt $customerIds.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $customerIds[$i]
}

# The second "MEREP" entry (row 15) was pasted in from a different source
# and keeps its own distinct (but visually identical) font.
$ws.Range("B15").Font.Name = "Calibri"
$ws.Range("B15").Font.Size = 11
$ws.Range("B15").Font.Color = 0

$ws.Range("B16").Select()
